$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

$ws.Range("D2:D3").Value = 0.0639
$ws.Range("E2:E3").Value = 0.0428
$ws.Range("G2:G3").Value = 0.1893732970027248
$ws.Range("H2:H3").Value = 0.1893732970027248
$ws.Range("I2:I3").Value = 0.215258855585831
$ws.Range("J2:J3").Value = 0.1778785752745191
$ws.Range("K2:K3").Value = 13.8
$ws.Range("L2:L3").Value = 0.1880108991825613
$ws.Range("M2:M3").Value = 11.4
$ws.Range("N2:N3").Value = 0.1011535048802129
$ws.Range("O2:O3").Value = 0.8260869565217391
$ws.Range("P2:P3").Value = 11.4
$ws.Range("Q2:Q3").Value = 0.1011535048802129
$ws.Range("R2:R3").Value = 0.8260869565217391
$ws.Range("U2:U3").Value = 7.46
$ws.Range("V2:V3").Value = 0.06619343389529725
$ws.Range("W2:W3").Value = 0.1073094867807154
$ws.Range("X2:X3").Value = 0.112048664881477
$ws.Range("Y2:Y3").Value = -0.004739178100761618
$ws.Range("Z2:Z3").Value = 0.6000654022236757
$ws.Range("AA2:AA3").Value = 0.1067387788190787
$ws.Range("AB2:AB3").Value = 0.112048664881477
$ws.Range("AC2:AC3").Value = -0.005309886062398364
$ws.Range("AG2:AG3").Value = -7.46
$ws.Range("AJ2:AJ3").Value = -0.07088559483086278
$ws.Range("AK2:AK3").Value = -0.05970866015687529
$ws.Range("AL2:AL3").Value = 1.09
$ws.Range("AM2:AM3").Value = 1.09
$ws.Range("AO2:AO3").Value = 14.4954128440367
$ws.Range("AP2:AP3").Value = -0.4691823899371069
$ws.Range("AQ2:AQ3").Value = 14.4954128440367
